$wb = $excel.ActiveWorkbook

# --- Work on the "Repayment schedule" sheet ---
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N ("Late"), shifting the
# existing N/O/P ("Late", "heading", "Outstanding") columns one to the
# right (O/P/Q).
$ws.Range("N1").EntireColumn.Insert()

# The newly inserted column keeps a custom (non bestFit) width of 11
# characters.
$ws.Range("N1").EntireColumn.ColumnWidth = 10.14

# Make "Repayment schedule" the active sheet/tab and set the new
# selection on it.
$ws.Activate()
$ws.Range("R8").Select()
